$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.112.67'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '2.302.72'
$ws.Range("E3").Value = '  +1.43%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '300.92'
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").Value = '97.76'
$ws.Range("E6").Value = '  +1.48%  '

$ws.Range("D7").Value = '0.508'
$ws.Range("E7").Value = '  +1.25%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '0.503'
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("D10").Value = '33.67'
$ws.Range("E10").Value = '  -0.12%  '

$ws.Range("E11").Value = '  +1.45%  '

$ws.Range("D12").Value = '49.37'
$ws.Range("E12").Value = '  -2.47%  '

$ws.Range("E13").Value = '  +3.10%  '

$ws.Range("D14").Value = '17.25'
$ws.Range("E14").Value = '  +13.41%  '

$ws.Range("D15").Value = '6.78'
$ws.Range("E15").Value = '  +2.17%  '

$ws.Range("D16").Value = '2.662.30'
$ws.Range("E16").Value = '  +1.54%  '

$ws.Range("D17").Value = '2.304.27'
$ws.Range("E17").Value = '  +1.62%  '

$ws.Range("D18").Value = '0.815'
$ws.Range("E18").Value = '  +4.11%  '

$ws.Range("D19").Value = '43.066.75'
$ws.Range("E19").Value = '  +1.79%  '

$ws.Range("D20").Value = '11.71'
$ws.Range("E20").Value = '  +2.71%  '

$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").Value = '6.08'
$ws.Range("E22").Value = '  +1.93%  '

$ws.Range("D23").Value = '67.92'
$ws.Range("E23").Value = '  +2.18%  '

$ws.Range("D24").Value = '237.08'
$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("D25").Value = '2.04'
$ws.Range("E25").Value = '  +5.70%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("D28").Value = '24.48'
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  +2.14%  '

$ws.Range("D30").Value = '167.64'
$ws.Range("E30").Value = '  +2.16%  '

$ws.Range("D31").Value = '34.18'
$ws.Range("E31").Value = '  +1.32%  '

$ws.Range("D32").Value = '9.14'
$ws.Range("E32").Value = '  +1.00%  '

$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("D34").Value = '4.97'
$ws.Range("E34").Value = '  +0.88%  '

$ws.Range("D35").Value = '4.65'
$ws.Range("E35").Value = '  +6.70%  '

$ws.Range("E36").Value = '  -1.01%  '

$ws.Range("D37").Value = '16.83'
$ws.Range("E37").Value = '  +4.95%  '

$ws.Range("D38").Value = '0.0701'
$ws.Range("E38").Value = '  +1.09%  '

$ws.Range("E39").Value = '  +0.99%  '

$ws.Range("D40").Value = '0.101'
$ws.Range("E40").Value = '  +1.51%  '

$ws.Range("D41").Value = '1.77'
$ws.Range("E41").Value = '  +0.70%  '

$ws.Range("E42").Value = '  +0.49%  '

$ws.Range("D43").Value = '2.37'
$ws.Range("E43").Value = '  -1.61%  '

$ws.Range("D44").Value = '1.980.05'
$ws.Range("E44").Value = '  +0.99%  '

$ws.Range("E45").Value = '  +0.81%  '

$ws.Range("D46").Value = '9.93'
$ws.Range("E46").Value = '  +2.59%  '

$ws.Range("D47").Value = '17.65'
$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("D48").Value = '2.86'
$ws.Range("E48").Value = '  +1.89%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.529.90'
$ws.Range("E49").Value = '  +1.02%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '53.50'
$ws.Range("E50").Value = '  +1.51%  '

$ws.Range("D51").Value = '4.60'
$ws.Range("E51").Value = '  -1.59%  '
